# Generate Report for Handoff
#
# Updates the localization-status workbook to reflect a new handoff run:
#   - b752bc2d-b7e0-4cd7-ada0-0b20e237869d is now "Ready for handoff"
#     (previously "In Translation"), with fresh handoff/target timestamps
#     and file names.
#   - 49ae636d-794c-44d6-bc6f-398bd2b0a69f moved into the "In Translation"
#     row/slot that b752bc2d... used to occupy.
#
# Helper: set a cell's value and (optionally) keep the matching hyperlink's
# displayed text in sync with it. Uses positional parameters throughout
# (named "-Param value" binding is unreliable in this host).
function Set-CellAndHyperlink($Sheet, $CellAddr, $NewValue, $HasHyperlink) {
    $Sheet.Range($CellAddr).Value = $NewValue
    if ($HasHyperlink) {
        $col = $CellAddr.Substring(0, 1)
        $row = $CellAddr.Substring(1)
        $target = ('$' + $col + '$' + $row)
        foreach ($hl in $Sheet.Hyperlinks) {
            if ($hl.Range.Address() -eq $target) {
                $hl.TextToDisplay = $NewValue
            }
        }
    }
}

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
Set-CellAndHyperlink $ws1 "A2" "b752bc2d-b7e0-4cd7-ada0-0b20e237869d.md" $true
Set-CellAndHyperlink $ws1 "A3" "49ae636d-794c-44d6-bc6f-398bd2b0a69f.md" $true
Set-CellAndHyperlink $ws1 "B3" "Ready for handoff" $false
Set-CellAndHyperlink $ws1 "C3" "Ready for handoff" $false
Set-CellAndHyperlink $ws1 "D3" "2016-18-19 14:18:24" $false

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
Set-CellAndHyperlink $ws2 "A2" "b752bc2d-b7e0-4cd7-ada0-0b20e237869d.md" $true
Set-CellAndHyperlink $ws2 "D2" "b752bc2d-b7e0-4cd7-ada0-0b20e237869d.6929df9fc004e14e92fa6ac8dc5da21f161e0e3d.zh-cn.xlf" $true
Set-CellAndHyperlink $ws2 "A3" "49ae636d-794c-44d6-bc6f-398bd2b0a69f.md" $true
Set-CellAndHyperlink $ws2 "C3" "Ready for handoff" $false
Set-CellAndHyperlink $ws2 "D3" "49ae636d-794c-44d6-bc6f-398bd2b0a69f.963d6a6cc6d1e9906195d87fe9b32672cf838d79.zh-cn.xlf" $true
Set-CellAndHyperlink $ws2 "E3" "2016-03-19 14:18:21" $false

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
Set-CellAndHyperlink $ws3 "A2" "b752bc2d-b7e0-4cd7-ada0-0b20e237869d.md" $true
Set-CellAndHyperlink $ws3 "D2" "b752bc2d-b7e0-4cd7-ada0-0b20e237869d.6929df9fc004e14e92fa6ac8dc5da21f161e0e3d.de-de.xlf" $true
Set-CellAndHyperlink $ws3 "A3" "49ae636d-794c-44d6-bc6f-398bd2b0a69f.md" $true
Set-CellAndHyperlink $ws3 "C3" "Ready for handoff" $false
Set-CellAndHyperlink $ws3 "D3" "49ae636d-794c-44d6-bc6f-398bd2b0a69f.963d6a6cc6d1e9906195d87fe9b32672cf838d79.de-de.xlf" $true
Set-CellAndHyperlink $ws3 "E3" "2016-03-19 14:18:24" $false

Write-Output "Generate Report for Handoff: applied."
